$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the "Production Percentage" header text (was "Diff. %")
$ws.Range("G8").Value = "Production Percentage"

# Widen column G to fit the new, longer header text
$ws.Columns.Item(7).ColumnWidth = 11.833333333333332

# Update the active selection to G8 (reflects where the user ended up after the edit)
$ws.Range("G8").Select()
